# Refresh the crypto price-tracker sheet with the latest scrape.
#
# The 'Price' (D) and 'Volume(1h)' (E) columns are stored as plain text
# (not numbers) in this workbook, including values that look numeric (e.g.
# "298.09") and percentages padded with spaces (e.g. "  -1.57%  "). Excel
# auto-converts a plain numeric-looking string typed into a General-format
# cell into a real number, so any replacement value that would parse as a
# number is written with a leading apostrophe to force text entry, matching
# how a person would type it into Excel to keep it as text. The leading
# apostrophe itself is a text-entry marker and is not stored in the cell.
#
# Rows 40/41 also swap places (Kaspa now ranks above Stellar).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '42.108.89'
$ws.Range('E2').Value = '  -1.57%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.255.55'
$ws.Range('E3').Value = '  -3.32%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.05%  '

# Row 5: BNB
$ws.Range('D5').Value = '''298.09'
$ws.Range('E5').Value = '  -2.72%  '

# Row 6: Solana
$ws.Range('D6').Value = '''94.22'
$ws.Range('E6').Value = '  -5.60%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.497'
$ws.Range('E7').Value = '  -2.24%  '

# Row 8: USDC
$ws.Range('E8').Value = '  +0.06%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  -3.48%  '

# Row 10: Avalanche
$ws.Range('D10').Value = '''32.82'
$ws.Range('E10').Value = '  -6.04%  '

# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.0784'
$ws.Range('E11').Value = '  -1.96%  '

# Row 12: OKB
$ws.Range('D12').Value = '''48.15'
$ws.Range('E12').Value = '  -7.53%  '

# Row 13: TRON
$ws.Range('E13').Value = '  +0.45%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  -2.69%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.604.33'
$ws.Range('E15').Value = '  -3.50%  '

# Row 16: Chainlink
$ws.Range('D16').Value = '''15.35'
$ws.Range('E16').Value = '  -2.80%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.259.75'
$ws.Range('E17').Value = '  -5.51%  '

# Row 18: Polygon
$ws.Range('D18').Value = '''0.773'
$ws.Range('E18').Value = '  -2.84%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '42.075.54'
$ws.Range('E19').Value = '  -1.52%  '

# Row 20: ShibaInu
$ws.Range('E20').Value = '  -2.09%  '

# Row 21: InternetComputer(DFINITY)
$ws.Range('E21').Value = '  -2.22%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '''5.99'
$ws.Range('E22').Value = '  -3.92%  '

# Row 23: Litecoin
$ws.Range('D23').Value = '''66.34'
$ws.Range('E23').Value = '  -2.00%  '

# Row 24: BitcoinCash
$ws.Range('D24').Value = '''232.59'
$ws.Range('E24').Value = '  -1.57%  '

# Row 25: ImmutableX
$ws.Range('D25').Value = '''1.92'
$ws.Range('E25').Value = '  -3.83%  '

# Row 26: Dai
$ws.Range('E26').Value = '  -0.02%  '

# Row 27: PancakeSwap
$ws.Range('D27').Value = '''2.45'
$ws.Range('E27').Value = '  -4.04%  '

# Row 28: EthereumClassic
$ws.Range('D28').Value = '''23.75'
$ws.Range('E28').Value = '  -5.11%  '

# Row 29: Monero
$ws.Range('D29').Value = '''166.14'
$ws.Range('E29').Value = '  +4.19%  '

# Row 30: Toncoin
$ws.Range('D30').Value = '''2.04'
$ws.Range('E30').Value = '  -12.14%  '

# Row 31: InjectiveProtocol
$ws.Range('E31').Value = '  -3.47%  '

# Row 32: Cosmos
$ws.Range('E32').Value = '  -3.64%  '

# Row 33: FirstDigitalUSD
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  -0.08%  '

# Row 34: Filecoin
$ws.Range('E34').Value = '  -3.60%  '

# Row 35: WEMIXToken
$ws.Range('E35').Value = '  -4.75%  '

# Row 36: Hedera
$ws.Range('E36').Value = '  -4.91%  '

# Row 37: RenderToken
$ws.Range('D37').Value = '''4.33'
$ws.Range('E37').Value = '  -5.61%  '

# Row 38: LidoDAOToken
$ws.Range('D38').Value = '''2.80'
$ws.Range('E38').Value = '  -5.36%  '

# Row 39: Celestia
$ws.Range('D39').Value = '''15.94'
$ws.Range('E39').Value = '  -8.09%  '

# Row 40: Kaspa
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.0981'
$ws.Range('E40').Value = '  -4.99%  '

# Row 41: Stellar
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '''0.109'
$ws.Range('E41').Value = '  -2.93%  '

# Row 42: ARBITRUM
$ws.Range('E42').Value = '  -8.35%  '

# Row 43: ApeXProtocol
$ws.Range('E43').Value = '  +2.95%  '

# Row 44: Maker
$ws.Range('D44').Value = '1.934.07'
$ws.Range('E44').Value = '  -4.41%  '

# Row 45: VeChain
$ws.Range('D45').Value = '''0.0278'
$ws.Range('E45').Value = '  -2.35%  '

# Row 46: EnergySwap
$ws.Range('D46').Value = '''17.33'
$ws.Range('E46').Value = '  -6.83%  '

# Row 47: FraxShare
$ws.Range('D47').Value = '''9.51'
$ws.Range('E47').Value = '  -7.49%  '

# Row 48: NEARProtocol
$ws.Range('E48').Value = '  -5.31%  '

# Row 49: HuobiToken
$ws.Range('E49').Value = '  -3.28%  '

# Row 50: RocketPoolETH
$ws.Range('D50').Value = '2.482.66'
$ws.Range('E50').Value = '  -2.87%  '

# Row 51: MultiversX
$ws.Range('D51').Value = '''51.99'
$ws.Range('E51').Value = '  -7.22%  '
